# ============================================================================
# Applies the "adding averages and more checks" update:
#  - Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
#    for every data row (the report was regenerated 8 days later: 16-Sep-2025
#    instead of 08-Sep-2025), flips rows 17 & 18 from VALID to NOT VALID and
#    gives them the same "not valid" pink highlight already used by row 19.
#  - Header / title fonts get a white font color.
#  - Exam Dashboard: narrows column E and rewrites the COMMENTS text.
# ============================================================================

$wb = $excel.ActiveWorkbook

$training = $wb.Worksheets.Item("Training Dashboard")
$exam     = $wb.Worksheets.Item("Exam Dashboard")

# ----------------------------------------------------------------------
# 1) Training Dashboard - per-row "PERIOD TO EXPIRE" / "LAST UPDATE" data
# ----------------------------------------------------------------------
$newH = @{
  3=426; 4=219; 5=447; 6=507; 7=483; 8=483; 9=525; 10=503; 11=489; 12=400;
  13=496; 14=399; 15=504; 16=489; 17=15; 18=26; 19=-181; 20=170; 21=155; 22=181
}

foreach ($r in $newH.Keys) {
    $training.Cells.Item($r, 8).Value2 = $newH[$r]

    # Column I holds a literal text date ("16-Sep-2025"), not a real date
    # value, so force Text formatting before assigning to stop Excel from
    # auto-converting the string into a date serial number.
    $icell = $training.Cells.Item($r, 9)
    $icell.NumberFormat = "@"
    $icell.Value2 = "16-Sep-2025"
}

# Rows 17 & 18 flip from VALID to NOT VALID
$training.Cells.Item(17, 10).Value2 = "NOT VALID"
$training.Cells.Item(18, 10).Value2 = "NOT VALID"

# Give rows 17 & 18 the same "NOT VALID" pink highlight formatting already
# used on row 19 (copy formats only, values are left untouched).
$training.Range("A19:K19").Copy() | Out-Null
$training.Range("A17:K17").PasteSpecial(-4122) | Out-Null
$training.Range("A18:K18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 2) Header / title font color -> white
#    (the title also drops its old explicit 14pt size so it now shares the
#    same bold/white font as the header row)
# ----------------------------------------------------------------------
foreach ($ws in @($training, $exam)) {
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count

    $ws.Range("A1").Font.Color = 16777215
    $ws.Range("A1").Font.Size = 11
    $ws.Range($ws.Cells.Item(2,1), $ws.Cells.Item(2,$lastCol)).Font.Color = 16777215
}

# ----------------------------------------------------------------------
# 3) Exam Dashboard - narrower comments column + updated remarks
# ----------------------------------------------------------------------
$exam.Columns.Item(5).ColumnWidth = 14.14

$exam.Cells.Item(3, 5).Value2 = "date is valid"
$exam.Cells.Item(4, 5).Value2 = "date is valid"
$exam.Cells.Item(5, 5).Value2 = "date is valid"
$exam.Cells.Item(6, 5).Value2 = "date is valid"
$exam.Cells.Item(7, 5).Value2 = "date is valid"
$exam.Cells.Item(8, 5).Value2 = "date is valid"

Write-Host "edit.ps1 completed"
